$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "50"
$ws.Range("B27").Value = "Crystal"
$ws.Range("A28").Value = "60"
$ws.Range("B28").Value = "Ability"

$ws.Range("A29").Select()
